# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit figures across the per-job "Leve Profits" sheets.
# (Mirrors an automated runner re-pulling Universalis prices.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35000
$ws.Range("J3").Value = 35000
$ws.Range("L3").Value = 35000
$ws.Range("N3").Value = -35228

$ws.Range("H13").Value = 5250
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5250
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5250
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5588

$ws.Range("H20").Value = 2919.4285
$ws.Range("I20").Value = 906
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 906
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -676
$ws.Range("N20").Value = -15460

$ws.Range("H33").Value = 691
$ws.Range("I33").Value = 783.4231
$ws.Range("K33").Value = 783.4231
$ws.Range("M33").Value = -554.4231

$ws.Range("H35").Value = 2919.4285
$ws.Range("I35").Value = 906
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 906
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = -527
$ws.Range("N35").Value = -15758

$ws.Range("H98").Value = 785.65
$ws.Range("I98").Value = 785.65
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 785.65
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 712.35
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 1856.591
$ws.Range("I100").Value = 1250.3846
$ws.Range("J100").Value = 2732.2222
$ws.Range("K100").Value = 1250.3846
$ws.Range("L100").Value = 2732.2222
$ws.Range("M100").Value = -709.3846000000001
$ws.Range("N100").Value = -3814.2222

$ws.Range("H102").Value = 35000
$ws.Range("J102").Value = 35000
$ws.Range("L102").Value = 35000
$ws.Range("N102").Value = -41490

$ws.Range("H122").Value = 785.65
$ws.Range("I122").Value = 785.65
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2356.95
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 93.05000000000018
$ws.Range("N122").ClearContents()

$ws.Range("H135").Value = 1396.76
$ws.Range("I135").Value = 1077.3889
$ws.Range("J135").Value = 2218
$ws.Range("K135").Value = 9696.500099999999
$ws.Range("L135").Value = 19962
$ws.Range("M135").Value = -7161.500099999999
$ws.Range("N135").Value = -25032

$ws.Range("H137").Value = 30562.795
$ws.Range("I137").Value = 784.92
$ws.Range("J137").Value = 113279.11
$ws.Range("K137").Value = 2354.76
$ws.Range("L137").Value = 339837.33
$ws.Range("M137").Value = 195.2400000000002
$ws.Range("N137").Value = -344937.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 31666.666
$ws.Range("J103").Value = 31666.666
$ws.Range("L103").Value = 31666.666
$ws.Range("N103").Value = -34010.666

$ws.Range("H122").Value = 3108.5806
$ws.Range("I122").Value = 2665.739
$ws.Range("J122").Value = 4381.75
$ws.Range("K122").Value = 7997.217000000001
$ws.Range("L122").Value = 13145.25
$ws.Range("M122").Value = -5547.217000000001
$ws.Range("N122").Value = -18045.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2341.7317
$ws.Range("I31").Value = 1650.2222
$ws.Range("J31").Value = 2882.913
$ws.Range("K31").Value = 1650.2222
$ws.Range("L31").Value = 2882.913
$ws.Range("M31").Value = -1355.2222
$ws.Range("N31").Value = -3472.913

$ws.Range("H34").Value = 2341.7317
$ws.Range("I34").Value = 1650.2222
$ws.Range("J34").Value = 2882.913
$ws.Range("K34").Value = 1650.2222
$ws.Range("L34").Value = 2882.913
$ws.Range("M34").Value = -1448.2222
$ws.Range("N34").Value = -3286.913

$ws.Range("H107").Value = 1333.8
$ws.Range("I107").Value = 1264.4445
$ws.Range("J107").Value = 1437.8334
$ws.Range("K107").Value = 1264.4445
$ws.Range("L107").Value = 1437.8334
$ws.Range("M107").Value = 655.5554999999999
$ws.Range("N107").Value = -5277.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 979.125
$ws.Range("J5").Value = 1520
$ws.Range("L5").Value = 4560
$ws.Range("N5").Value = -4784

$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868

$ws.Range("H107").Value = 96.695656
$ws.Range("I107").Value = 95.30768999999999
$ws.Range("J107").Value = 98.5
$ws.Range("K107").Value = 285.92307
$ws.Range("L107").Value = 295.5
$ws.Range("M107").Value = 1634.07693
$ws.Range("N107").Value = -4135.5

$ws.Range("H127").Value = 922.44446
$ws.Range("J127").Value = 922.44446
$ws.Range("L127").Value = 2767.33338
$ws.Range("N127").Value = -12687.33338

$ws.Range("H135").Value = 979.125
$ws.Range("J135").Value = 1520
$ws.Range("L135").Value = 13680
$ws.Range("N135").Value = -18750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 925.7143
$ws.Range("I13").Value = 196
$ws.Range("J13").Value = 2750
$ws.Range("K13").Value = 196
$ws.Range("L13").Value = 2750
$ws.Range("M13").Value = -57
$ws.Range("N13").Value = -3028

$ws.Range("H52").Value = 24946
$ws.Range("I52").Value = 2500
$ws.Range("J52").Value = 25404.082
$ws.Range("K52").Value = 2500
$ws.Range("L52").Value = 25404.082
$ws.Range("M52").Value = -2241
$ws.Range("N52").Value = -25922.082

$ws.Range("H105").Value = 34998
$ws.Range("J105").Value = 34998
$ws.Range("L105").Value = 34998
$ws.Range("N105").Value = -41986

$ws.Range("H122").Value = 4251
$ws.Range("I122").Value = 4251
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12753
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10303
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2074.6316
$ws.Range("I7").Value = 1257.8
$ws.Range("J7").Value = 2982.2222
$ws.Range("K7").Value = 1257.8
$ws.Range("L7").Value = 2982.2222
$ws.Range("M7").Value = -1145.8
$ws.Range("N7").Value = -3206.2222

$ws.Range("H46").Value = 3700
$ws.Range("I46").Value = 3700
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3700
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3512
$ws.Range("N46").ClearContents()

$ws.Range("H106").Value = 21668.857
$ws.Range("J106").Value = 23613.834
$ws.Range("L106").Value = 23613.834
$ws.Range("N106").Value = -26137.834

$ws.Range("H122").Value = 5879.35
$ws.Range("I122").Value = 6306.533
$ws.Range("J122").Value = 4597.8
$ws.Range("K122").Value = 18919.599
$ws.Range("L122").Value = 13793.4
$ws.Range("M122").Value = -16469.599
$ws.Range("N122").Value = -18693.4

$ws.Range("H126").Value = 2074.6316
$ws.Range("I126").Value = 1257.8
$ws.Range("J126").Value = 2982.2222
$ws.Range("K126").Value = 3773.4
$ws.Range("L126").Value = 8946.6666
$ws.Range("M126").Value = -1303.4
$ws.Range("N126").Value = -13886.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 25320
$ws.Range("J63").Value = 25320
$ws.Range("L63").Value = 25320
$ws.Range("N63").Value = -26568

$ws.Range("H66").Value = 25320
$ws.Range("J66").Value = 25320
$ws.Range("L66").Value = 75960
$ws.Range("N66").Value = -82200

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H100").Value = 30830.1
$ws.Range("I100").Value = 120000.4
$ws.Range("J100").Value = 1106.6666
$ws.Range("K100").Value = 240000.8
$ws.Range("L100").Value = 2213.3332
$ws.Range("M100").Value = -239459.8
$ws.Range("N100").Value = -3295.3332
